$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 5 (Honeycomb / VP EMEA / Laura Cleaver / 1st Interview),
# shifting row 6 (Blockaid / Head of EMEA / Sam Wallis / 1st Interview) up to row 5.
$ws.Rows.Item(5).Delete()

# Update the status of the now-row-5 (Blockaid) entry to "2nd Interview".
$ws.Range("E5").Value = "2nd Interview"
